# Adds the new "tọa thủ / đồng cung ... tại cung đối Phụ Mẫu" entries
# to Sheet2, continuing the existing A/B mirrored-string table.
# New data starts at row 4329 (row 4328 is intentionally left absent,
# matching the pre-existing row-numbering gaps already present in the sheet)
# and runs through row 4433 (105 new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    'Tử Vi tọa thủ tại cung đối Phụ Mẫu',
    'Thiên Cơ tọa thủ tại cung đối Phụ Mẫu',
    'Thái Dương tọa thủ tại cung đối Phụ Mẫu',
    'Vũ Khúc tọa thủ tại cung đối Phụ Mẫu',
    'Thiên Đồng tọa thủ tại cung đối Phụ Mẫu',
    'Liêm Trinh tọa thủ tại cung đối Phụ Mẫu',
    'Thiên Phủ tọa thủ tại cung đối Phụ Mẫu',
    'Thái Âm tọa thủ tại cung đối Phụ Mẫu',
    'Tham Lang tọa thủ tại cung đối Phụ Mẫu',
    'Cự Môn tọa thủ tại cung đối Phụ Mẫu',
    'Thiên Tướng tọa thủ tại cung đối Phụ Mẫu',
    'Thiên Lương tọa thủ tại cung đối Phụ Mẫu',
    'Thất Sát tọa thủ tại cung đối Phụ Mẫu',
    'Phá Quân tọa thủ tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thiên Cơ tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thái Dương tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Vũ Khúc tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thiên Đồng tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Liêm Trinh tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Tử Vi đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thái Dương tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Vũ Khúc tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thiên Đồng tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Liêm Trinh tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thiên Cơ đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Vũ Khúc tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thiên Đồng tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Liêm Trinh tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thái Dương đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thiên Đồng tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Liêm Trinh tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Vũ Khúc đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Liêm Trinh tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thiên Đồng đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Thiên Phủ tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Liêm Trinh đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Thái Âm tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thiên Phủ đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Tham Lang tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thái Âm đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Tham Lang đồng cung Cự Môn tại cung đối Phụ Mẫu',
    'Tham Lang đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Tham Lang đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Tham Lang đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Tham Lang đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Cự Môn đồng cung Thiên Tướng tại cung đối Phụ Mẫu',
    'Cự Môn đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Cự Môn đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Cự Môn đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thiên Tướng đồng cung Thiên Lương tại cung đối Phụ Mẫu',
    'Thiên Tướng đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thiên Tướng đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thiên Lương đồng cung Thất Sát tại cung đối Phụ Mẫu',
    'Thiên Lương đồng cung Phá Quân tại cung đối Phụ Mẫu',
    'Thất Sát đồng cung Phá Quân tại cung đối Phụ Mẫu'
)

$startRow = 4329
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $text = $newValues[$i]
    $ws.Cells.Item($row, 1).Value = $text
    $ws.Cells.Item($row, 2).Value = $text
}

# Reflect the new extent in the sheet's selection/active range, mirroring
# how Excel leaves the selection after filling the new block: the 14
# single-star "tọa thủ" rows (4329-4342) were entered first, then the 91
# paired-star "đồng cung" rows (4343-4433) were filled/pasted as one block,
# which is what's left selected.
$lastRow = $startRow + $newValues.Count - 1
$pairBlockStart = $startRow + 14
[void]$ws.Range("B" + $pairBlockStart + ":B" + $lastRow).Select()
